# data : case 1
# Update the "case" data row (row 1) with the new sample values and
# resize the columns to match the new case's layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 1 cell values -------------------------------------------------
$ws.Range("B1").Value = 3
$ws.Range("C1").Value = 11
$ws.Range("D1").Value = 8
$ws.Range("E1").Value = 6
$ws.Range("F1").Value = 2
$ws.Range("G1").Value = 19
$ws.Range("H1").Value = 11
$ws.Range("I1").Value = 0.09501
$ws.Range("J1").Value = 0.085
$ws.Range("K1").Value = 0.099

# --- Resize columns to match the new case layout ------------------------------
# (Column widths snap to the workbook's character-width grid, same as Excel's
# own ColumnWidth dialog, so the values below are chosen to land on the
# closest grid point to the target widths.)
$ws.Range("D1:F1").ColumnWidth = 1.3333333333333333
$ws.Range("H1").ColumnWidth = 2.3333333333333335
$ws.Range("I1").ColumnWidth = 6.833333333333333
$ws.Range("J1:K1").ColumnWidth = 4.833333333333333
